# Applies the "added my own version of random walks" commit:
#  - ChosenTreatment!A2/B2 get new condition/treatment dictionaries
#  - Summary!A2 (NumSubgroups) goes from 21 to 26
#  - Subgroups sheet gets a refreshed set of 26 subgroup rows (was 21)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ChosenTreatment sheet
# ---------------------------------------------------------------------
$wsChosen = $wb.Worksheets.Item("ChosenTreatment")
$wsChosen.Range("A2").Value = "{'SexualOrientation': 'Straight or heterosexual'}"
$wsChosen.Range("B2").Value = "{'DevType': 'Back-end developer'}"

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A2").Value = 26

# ---------------------------------------------------------------------
# Subgroups sheet
# ---------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("Subgroups")

$rows = @(
    @("{'Hobby': np.int64(1)}", 28842, 5516.21828962618, -1197.137411587141),
    @("{'Student': np.int64(1)}", 29526, 5617.758312088003, -1095.597389125319),
    @("{'UndergradMajor': np.int64(2)}", 22173, 3542.149042873403, -3171.206658339918),
    @("{'Gender': np.int64(1)}", 33253, 5920.984439082945, -792.3712621303766),
    @("{'RaceEthnicity': np.int64(1)}", 27379, 6509.936391650183, -203.4193095631381),
    @("{'Dependents': np.int64(2)}", 24167, 7072.162954162452, 358.8072529491301),
    @("{'HDI': np.int64(1)}", 27581, 8280.034330059932, 1566.678628846611),
    @("{'Hobby': np.int64(1), 'Student': np.int64(1)}", 23894, 4208.446115890501, -2504.90958532282),
    @("{'Hobby': np.int64(1), 'Gender': np.int64(1)}", 27456, 5084.527001737077, -1628.828699476245),
    @("{'Hobby': np.int64(1), 'RaceEthnicity': np.int64(1)}", 22529, 5345.204113474854, -1368.151587738467),
    @("{'Hobby': np.int64(1), 'HDI': np.int64(1)}", 22660, 7025.684796894849, 312.3290956815272),
    @("{'Student': np.int64(1), 'Gender': np.int64(1)}", 27820, 4860.761884746209, -1852.593816467112),
    @("{'Student': np.int64(1), 'RaceEthnicity': np.int64(1)}", 23264, 5618.412897795355, -1094.942803417966),
    @("{'Student': np.int64(1), 'HDI': np.int64(1)}", 23773, 7095.77558069452, 382.4198794811982),
    @("{'UndergradMajor': np.int64(2), 'Gender': np.int64(1)}", 21070, 2140.616974317039, -4572.738726896283),
    @("{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1)}", 25910, 6471.012774252379, -242.3429269609423),
    @("{'Gender': np.int64(1), 'Dependents': np.int64(2)}", 22568, 6660.009269424866, -53.34643178845545),
    @("{'Gender': np.int64(1), 'HDI': np.int64(1)}", 25944, 7219.315451383841, 505.9597501705193),
    @("{'RaceEthnicity': np.int64(1), 'HDI': np.int64(1)}", 23146, 7669.000695259265, 955.6449940459433),
    @("{'Hobby': np.int64(1), 'Student': np.int64(1), 'Gender': np.int64(1)}", 22758, 4081.320153739394, -2632.035547473928),
    @("{'Hobby': np.int64(1), 'Gender': np.int64(1), 'RaceEthnicity': np.int64(1)}", 21544, 5551.102665918384, -1162.253035294938),
    @("{'Hobby': np.int64(1), 'Gender': np.int64(1), 'HDI': np.int64(1)}", 21574, 6412.706957617019, -300.6487435963027),
    @("{'Student': np.int64(1), 'Gender': np.int64(1), 'RaceEthnicity': np.int64(1)}", 22014, 5982.017555487582, -731.3381457257392),
    @("{'Student': np.int64(1), 'Gender': np.int64(1), 'HDI': np.int64(1)}", 22370, 6176.689060932703, -536.6666402806186),
    @("{'Student': np.int64(1), 'RaceEthnicity': np.int64(1), 'HDI': np.int64(1)}", 20128, 6694.06933007496, -19.28637113836157),
    @("{'Gender': np.int64(1), 'RaceEthnicity': np.int64(1), 'HDI': np.int64(1)}", 21900, 7566.667087338407, 853.3113861250858)
)

$r = 2
foreach ($row in $rows) {
    $wsSub.Cells.Item($r, 1).Value = $row[0]
    $wsSub.Cells.Item($r, 2).Value = $row[1]
    $wsSub.Cells.Item($r, 3).Value = $row[2]
    $wsSub.Cells.Item($r, 4).Value = $row[3]
    $r++
}
